# Set the "Industries" column (H) values to 0 for rows 25 through 176.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 25; $row -le 176; $row++) {
    $ws.Range("H$row").Value = 0
}
